$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted before the existing row 12,
# pushing every subsequent record down by one row (row 12 -> 13, ...,
# row 56 -> 57). The sheet's used range grows from A1:R56 to A1:R57.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with this week's record.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Vega Modelo de Temuco"
$ws.Range("C12").Value = "La Araucanía"
$ws.Range("D12").Value = 44910
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 100112042
$ws.Range("G12").Value = "Locoto"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 180
$ws.Range("K12").Value = 2500
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2500
$ws.Range("N12").Value = "$/kilo"
$ws.Range("O12").Value = "Región de Arica y Parinacota"
$ws.Range("P12").Value = 2500
$ws.Range("Q12").Value = 1
$ws.Range("R12").Value = "Hortaliza"
